# Updated test data for 5,24,40V,BatteryStandby and AC Calculations test cases
#
# The only functional change is on the "Add Panels" sheet: cell B4 (which was
# an empty, bordered cell) now holds the new "Loading Details Name" value
# "NGC-601/T1461OR TC-212" (a new shared-string entry), and the cell loses
# its previous border formatting (no style index on save). The active
# selection in the sheet also moves to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing border/style from B4 before writing the new text so the
# saved cell carries no explicit style index (matches the target workbook).
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "NGC-601/T1461OR TC-212"

# Reflect the new active cell/selection recorded in the sheet view.
$ws.Range("B4").Select()
